$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.702.01'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.51%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.770.69'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +5.51%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '116.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '333.80'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.10%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.576'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.08'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0863'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.53%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.30'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.87%  '
$ws.Range("E13").Value = '  +2.24%  '
$ws.Range("E14").Value = '  +5.30%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.205.14'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.61%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.768.81'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.71%  '
$ws.Range("E17").Value = '  +4.20%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.659.64'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.65%  '
$ws.Range("E19").Value = '  +13.73%  '
$ws.Range("E20").Value = '  +5.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.86'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.91%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0977'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '278.72'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.71%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.90'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.38%  '
$ws.Range("E25").Value = '  +6.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.84'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.47%  '
$ws.Range("E27").Value = '  -0.12%  '
$ws.Range("E28").Value = '  -0.69%  '
$ws.Range("E29").Value = '  +0.26%  '
$ws.Range("E30").Value = '  +3.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.17'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '50.28'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.76%  '
$ws.Range("E34").Value = '  +0.94%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '19.31'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.61%  '
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("E37").Value = '  +3.15%  '
$ws.Range("E38").Value = '  +2.99%  '
$ws.Range("E39").Value = '  +4.81%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0356'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +9.74%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '128.55'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '23.43'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.31'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +7.83%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.114'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.62%  '
$ws.Range("E45").Value = '  +17.81%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.090.18'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.08%  '
$ws.Range("E47").Value = '  +3.54%  '
$ws.Range("E48").Value = '  +3.91%  '
$ws.Range("E49").Value = '  +7.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '60.42'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.00%  '
$ws.Range("E51").Value = '  -0.36%  '
